# [tsomsomm] : add CharacterData
# Replaces the two MapNpcData rows (2-3) with new CharacterData /
# CutsceneData localisation keys, and appends the remaining new
# CutsceneData rows (4-10), growing the used range from A1:D3 to A1:D10.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LocalData")

$rows = @(
    @("CharacterData.CharacterInfo.1000003",    "테스트용 캐릭터"),
    @("CutsceneData.Cutscene.1000001.1",        "ShowDialog"),
    @("CutsceneData.Cutscene.1000001.2",        "ShowDialogFlipped"),
    @("CutsceneData.Cutscene.1000002.4",        "ShowDialog"),
    @("CutsceneData.Cutscene.1000002.5",        "ShowDialogFlipped"),
    @("CutsceneData.CutsceneInfo.1000001",      "컷신이름1"),
    @("CutsceneData.CutsceneInfo.1000002",      "컷신이름2"),
    @("CutsceneData.CutsceneSelection.1000004", "선택지1"),
    @("CutsceneData.CutsceneSelection.1000005", "선택지2")
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $r++
}
